$p = $ppt.ActivePresentation

# Layout 2 = "Title and Content" (Title placeholder + single Content placeholder)
$ppLayoutText = 2

$slidesData = @(
    @{
        Title = "Introduction"
        Bullets = @(
            "- Introduction of the problem of Generalized Category Discovery (GCD) in fine-grained datasets.",
            "- Description of the challenges faced in fine-grained classification tasks due to high class similarities and intra-class variances."
        )
    },
    @{
        Title = "Related Works"
        Bullets = @(
            "- Overview of Novel Category Discovery (NCD) and its relevance to XCon's objectives.",
            "- Discussion on previous approaches utilizing transfer learning, self-supervision, and contrastive learning in categorizing unseen classes."
        )
    },
    @{
        Title = "Methodology"
        Bullets = @(
            "- Explanation of Expert-Contrastive Learning (XCon) approach.",
            "- Description of partitioning the dataset into expert sub-datasets using k-means clustering.",
            "- Details on supervised and unsupervised contrastive learning across full and sub-datasets."
        )
    },
    @{
        Title = "Experiments and Results"
        Bullets = @(
            "- Evaluation of XCon on CIFAR-10/100, ImageNet-100, CUB-200, Stanford Cars, FGVC-Aircraft, and Oxford-IIIT Pet datasets.",
            "- Showcase of state-of-the-art performance in fine-grained category discovery benchmarks.",
            "- Implementation details using ViT-B-16 model and DINO-pretrained parameters."
        )
    },
    @{
        Title = "Ablation Studies"
        Bullets = @(
            "- Overview of ablation studies conducted to validate components of XCon.",
            "- Discussion on the impact of weight of fine-grained loss and the number of sub-datasets.",
            "- Insights into optimal settings for different datasets based on experimental results."
        )
    },
    @{
        Title = "Conclusion"
        Bullets = @(
            "- Summary of XCon's significance in fine-grained category discovery.",
            "- Mention of the method's success across various benchmarks and its contribution to unsupervised and semi-supervised learning in fine-grained classification.",
            "- Encouragement for further exploration and adaptation of XCon within the research community."
        )
    }
)

$index = 1
foreach ($data in $slidesData) {
    $slide = $p.Slides.Add($index, $ppLayoutText)

    $titleShape = $slide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Text = $data.Title

    $bodyShape = $slide.Shapes.Item(2)
    $bodyShape.TextFrame.TextRange.Text = [string]::Join("`r", $data.Bullets)

    $index = $index + 1
}
